$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B79 was stored as a text/inline string "4" - convert it to a real number 4
$ws.Range("B79").Value = 4

# Add new row 80 with the new annotation entry
$ws.Range("A80").Value = "Ruilin"

# B80 must stay textual ("3"), even though it looks numeric. Force the cell
# to a Text format before assigning so Excel keeps it as a string, then
# reset the style back to Normal so no stray number-format styling remains.
$b80 = $ws.Range("B80")
$b80.NumberFormat = "@"
$b80.Value = "3"
$b80.Style = "Normal"

$ws.Range("C80").Value = "无"
$ws.Range("D80").Value = "DIS"
$ws.Range("E80").Value = "RES"
$ws.Range("F80").Value = "10db20b3-10d6-4fc8-9ec8-790c85de9594"
$ws.Range("G80").Value = "ByQZjx-0-_annotated.xlsx"
$ws.Range("H80").Value = "We suspect there will be no longer improvement, as we observed convergence in the controller's samples small entropy."
